$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.996.01'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '3.092.47'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.82'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.52'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.085.48'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  +6.62%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.61'
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').Value = '  +2.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.05'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000217'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '3.586.07'
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('D16').Value = '63.029.65'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '3.083.37'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '502.80'
$ws.Range('E19').Value = '  +3.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.50'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.704'
$ws.Range('E22').Value = '  +4.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.24'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.28'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.30'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.75'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.19'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.96'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.21'
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.51'
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.24'
$ws.Range('E34').Value = '  +14.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '526.83'
$ws.Range('E35').Value = '  -7.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.90'
$ws.Range('E36').Value = '  +2.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.18'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0413'
$ws.Range('E38').Value = '  +4.07%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.121'
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0788'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.038.44'
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.06'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.63'
$ws.Range('E43').Value = '  -5.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.253'
$ws.Range('E44').Value = '  +6.09%  '
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.95'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('B48').Value = 'CoreDAO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.44'
$ws.Range('E48').Value = '  +76.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.79'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('B51').Value = 'PEPE'
$ws.Range('C51').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D51').Value = '0.0₃0501'
$ws.Range('E51').Value = '  -1.74%  '
